$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32: new clue cell + puzzle grid values
$ws.Range("A32").Value = "1,1"
$ws.Range("B32").Value = 1
$ws.Range("D32").Value = 1

# Row 33: new grid values
$ws.Range("B33").Value = 1
$ws.Range("E33").Value = 1

# Row 34: new grid values
$ws.Range("B34").Value = 1
$ws.Range("F34").Value = 1

# Row 35 (new)
$ws.Range("C35").Value = 1
$ws.Range("E35").Value = 1

# Row 36 (new)
$ws.Range("C36").Value = 1
$ws.Range("F36").Value = 1

# Row 37 (new)
$ws.Range("D37").Value = 1
$ws.Range("F37").Value = 1

# Row 38 (new)
$ws.Range("I38").Value = 1
$ws.Range("P38").Value = 1

# Row 39 (new) - combinator output row 1
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = 1
$ws.Range("S39").Value = 0
$ws.Range("T39").Value = 0

# Row 40 (new) - combinator output row 2
$ws.Range("H40").Value = 1
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0
$ws.Range("P40").Value = 0
$ws.Range("Q40").Value = 0
$ws.Range("R40").Value = 0
$ws.Range("S40").Value = 1
$ws.Range("T40").Value = 0

# Row 41 (new) - combinator output row 3
$ws.Range("H41").Value = 2
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("S41").Value = 0
$ws.Range("T41").Value = 1

# Row 43 (new) - row 42 intentionally left blank/untouched
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0

# Row 44 (new)
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1

# Row 45 (new)
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 2

# Row 46 (new)
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = 1

# Row 47 (new)
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = 2

# Row 48 (new)
$ws.Range("I48").Value = 2
$ws.Range("J48").Value = 2

# Selection moves to K43
$ws.Range("K43").Select()
